$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.127.24"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "1.629.59"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.515"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "1.614.94"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.540"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.37%  "
$ws.Range("D16").Value = "27.095.09"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").Value = "0.0₃0732"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("E18").Value = "  -2.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0503"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").Value = "1.315.43"
$ws.Range("E33").Value = "  +3.81%  "
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.535"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.41%  "
$ws.Range("D43").Value = "1.766.30"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("E47").Value = "  +16.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.789"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +16.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0515"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.67%  "
$ws.Range("E51").Value = "  +0.09%  "
